# Replace the sample "intake" sheet's data with the new No/Age/Sallary
# unit-test dataset, apply an integer number format to the No/Age columns,
# move the sheet selection, and set the page setup (paper size / orientation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intake")

# Wipe the old Packet/Product/wt/color sample data (A1:D3).
$ws.Range("A1:D3").ClearContents()

# New headers.
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Sallary"

# New sample rows.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 60000

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 40000

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 3500

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 60000

# No / Age columns get an integer ("0") number format.
$ws.Range("A1:B5").NumberFormat = "0"

# Move the active selection.
$ws.Range("H10").Select()

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
